$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Mmp1a -> Itga2, FAPs -> ECs) updated TPM-derived values
$ws.Range("M2").Value = 3.425446666666666
$ws.Range("N2").Value = 10.27634
$ws.Range("O2").Value = 0.6657953389778073
$ws.Range("P2").Value = 0.6657953389778073
$ws.Range("Q2").Value = 0.2113295066533333
$ws.Range("R2").Value = 1.90196555988
$ws.Range("S2").Value = 0.6657953389778073
$ws.Range("T2").Value = 0.6657953389778073

# Row 3 (FAPs -> FAPs)
$ws.Range("O3").Value = 0.2094791321596951
$ws.Range("P3").Value = 0.2094791321596952
$ws.Range("S3").Value = 0.2094791321596951
$ws.Range("T3").Value = 0.2094791321596952

# Row 4 (FAPs -> MuSCs)
$ws.Range("M4").Value = 0.62317
$ws.Range("N4").Value = 1.86951
$ws.Range("O4").Value = 0.1211239647746572
$ws.Range("P4").Value = 0.1211239647746572
$ws.Range("Q4").Value = 0.03844584998
$ws.Range("R4").Value = 0.34601264982
$ws.Range("S4").Value = 0.1211239647746572
$ws.Range("T4").Value = 0.1211239647746572

# Row 5 (FAPs -> Resolving-Mac)
$ws.Range("M5").Value = 0.01852966666666667
$ws.Range("N5").Value = 0.055589
$ws.Range("O5").Value = 0.003601564087840353
$ws.Range("P5").Value = 0.003601564087840353
$ws.Range("Q5").Value = 0.001143169255333333
$ws.Range("R5").Value = 0.010288523298
$ws.Range("S5").Value = 0.003601564087840353
$ws.Range("T5").Value = 0.003601564087840353
